$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new Youtube links in column D for rows 5-8
$ws.Range("D5").Value = "https://youtu.be/WD5Un8e3EXQ"
$ws.Range("D6").Value = "https://youtu.be/gg1fX1t5OkA"
$ws.Range("D7").Value = "https://youtu.be/or4SKW0pwBY"
$ws.Range("D8").Value = "https://youtu.be/1fSfcE1hhzI"

# Update the selected/active cell to D10 (as reflected in the saved view state)
$ws.Range("D10").Select()
